$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Merge the two rows describing the "2nd Look" retrolective MRP re-evaluation
# option into a single combined label ...
$ws.Range("C181").Value = "1 - 2nd Look / Zweite MRP-Bewertung durchführen"

# ... and drop the now-redundant row that used to hold the German-only label
# (everything below shifts up by one row).
$ws.Rows("182:182").Delete()
